$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '68.727.19'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.82%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.830.77'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.77%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.71'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.45%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '169.66'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.825.95'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.64%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.523'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.161'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.42'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.01%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.453'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000261'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -2.67%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '36.53'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.64%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.479.57'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.827.56'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.812.62'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.88%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '18.11'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.67%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.02'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.64%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.95'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '468.22'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.705'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.14%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000154'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +8.79%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.49'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.42%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.19'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -2.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.00'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.97%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.20'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +0.31%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.78'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.37'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -0.51%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '30.29'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '9.20'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.92%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.785.49'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.999'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.101'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.77%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.47'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.43%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.73%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.83'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.07%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('B44').Value = 'TheGraph'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.304'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.76%  '
$ws.Range('B45').Value = 'Arweave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '44.11'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +17.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.94'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.30%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '46.85'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.96%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.46'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '146.47'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '392.58'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.805.50'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +4.45%  '
